$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5986
$ws.Range("F5").Value = 3044
$ws.Range("F6").Value = 450
$ws.Range("F8").Value = 60
$ws.Range("F9").Value = 47
$ws.Range("F11").Value = 354
$ws.Range("F12").Value = 4531
$ws.Range("F13").Value = 4531
$ws.Range("F15").Value = 104
$ws.Range("F16").Value = 139
$ws.Range("F17").Value = 25
$ws.Range("F19").Value = 88
$ws.Range("F20").Value = 7081
$ws.Range("F23").Value = 293
$ws.Range("F24").Value = 488
$ws.Range("F25").Value = 1290
$ws.Range("F27").Value = 1662
$ws.Range("F30").Value = 6069
$ws.Range("F31").Value = 127
$ws.Range("F33").Value = 106
$ws.Range("F34").Value = 89
$ws.Range("F36").Value = 6186
$ws.Range("F38").Value = 194
$ws.Range("F41").Value = 12
$ws.Range("F42").Value = 2430
$ws.Range("F43").Value = 27
$ws.Range("F44").Value = 49
$ws.Range("F45").Value = 1013
$ws.Range("F46").Value = 24
$ws.Range("F47").Value = 381
$ws.Range("F48").Value = 2091
$ws.Range("F49").Value = 25

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 217
$ws.Range("F9").Value = 43

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1431

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1431
$ws.Range("F3").Value = 5986
$ws.Range("F4").Value = 5986
$ws.Range("F5").Value = 3044
$ws.Range("F6").Value = 450
$ws.Range("F7").Value = 60
$ws.Range("F8").Value = 217
$ws.Range("F9").Value = 47
$ws.Range("F11").Value = 354
$ws.Range("F12").Value = 4531
$ws.Range("F13").Value = 4531
$ws.Range("F15").Value = 104
$ws.Range("F16").Value = 139
$ws.Range("F17").Value = 25
$ws.Range("F19").Value = 88
$ws.Range("F20").Value = 7081
$ws.Range("F23").Value = 489
$ws.Range("F24").Value = 1290
$ws.Range("F27").Value = 1662
$ws.Range("F30").Value = 43
$ws.Range("F31").Value = 6069
$ws.Range("F32").Value = 127
$ws.Range("F35").Value = 106
$ws.Range("F36").Value = 89
$ws.Range("F38").Value = 6186
$ws.Range("F40").Value = 194
$ws.Range("F42").Value = 12
$ws.Range("F44").Value = 2430
$ws.Range("F45").Value = 27
$ws.Range("F46").Value = 1013
$ws.Range("F47").Value = 24
$ws.Range("F48").Value = 381
$ws.Range("F49").Value = 2091
$ws.Range("F50").Value = 25
